# Adds a new weekly batch of "Sandia" (watermelon) price records for
# Mercado Mayorista Lo Valledor de Santiago, inserted right before the
# existing 2021-12-22 block, and shifts all subsequent rows down.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 blank rows at 977 - everything currently at/after row 977
# (old rows 977:991) moves down to 981:995.
$ws.Rows("977:980").Insert()

# Common (unchanging) column values for every data row in this sheet.
$mercadoId   = 6
$mercado     = "Mercado Mayorista Lo Valledor de Santiago"
$region      = "Metropolitana"
$codreg      = 13
$categoriaId = 100112028
$categoria   = "Sandia"
$variedad    = "Sin especificar"
$unidadCom   = '$/unidad'
$kgOUnidades = 1
$clasif      = "Hortaliza"

$fecha = 44628

$newRows = @(
    @{ Row = 977; Calidad = "Extra";   Volumen = 3200; PMin = 2800; PMax = 3000; PProm = 2881; Origen = "Región Metropolitana" },
    @{ Row = 978; Calidad = "Primera"; Volumen = 3800; PMin = 2300; PMax = 2500; PProm = 2384; Origen = "Región Metropolitana" },
    @{ Row = 979; Calidad = "Segunda"; Volumen = 2500; PMin = 1700; PMax = 2000; PProm = 1820; Origen = "Región Metropolitana" },
    @{ Row = 980; Calidad = "Tercera"; Volumen = 2700; PMin = 1200; PMax = 1300; PProm = 1241; Origen = "Región Metropolitana" }
)

foreach ($r in $newRows) {
    $arr = New-Object 'object[,]' 1,18
    $arr[0,0]  = $mercadoId
    $arr[0,1]  = $mercado
    $arr[0,2]  = $region
    $arr[0,3]  = $fecha
    $arr[0,4]  = $codreg
    $arr[0,5]  = $categoriaId
    $arr[0,6]  = $categoria
    $arr[0,7]  = $variedad
    $arr[0,8]  = $r.Calidad
    $arr[0,9]  = $r.Volumen
    $arr[0,10] = $r.PMin
    $arr[0,11] = $r.PMax
    $arr[0,12] = $r.PProm
    $arr[0,13] = $unidadCom
    $arr[0,14] = $r.Origen
    $arr[0,15] = $r.PProm
    $arr[0,16] = $kgOUnidades
    $arr[0,17] = $clasif

    $rowNum = $r.Row
    $ws.Range("A$($rowNum):R$($rowNum)").Value = $arr
}
